# Update "Pais" (COVID countries) sheet with refreshed figures and re-sorted rows.
# Source: paises.xlsx diff - "Update countries & provincias Spain"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh "last updated" timestamp
$ws.Range("A1").Value = 'Datos actualizados a 11 de Agosto de 2020 a las 18:05'

# Row 4: refreshed figures
$ws.Range("B4").Value = 5263827
$ws.Range("C4").Value = 12381
$ws.Range("D4").Value = 2717778
$ws.Range("E4").Value = 2379342
$ws.Range("G4").Value = 515
$ws.Range("H4").Value = 166707

# Row 6: refreshed figures
$ws.Range("B6").Value = 2322755
$ws.Range("C6").Value = 55602
$ws.Range("D6").Value = 1633356
$ws.Range("E6").Value = 643337
$ws.Range("G6").Value = 709
$ws.Range("H6").Value = 46062

# Row 12: refreshed figures
$ws.Range("B12").Value = 376616
$ws.Range("C12").Value = 1572
$ws.Range("D12").Value = 349541
$ws.Range("E12").Value = 16897
$ws.Range("G12").Value = 39
$ws.Range("H12").Value = 10178

# Row 20: refreshed figures
$ws.Range("B20").Value = 251237
$ws.Range("C20").Value = 412
$ws.Range("D20").Value = 202461
$ws.Range("E20").Value = 13561
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 35215

# Row 28: refreshed figures
$ws.Range("B28").Value = 113646
$ws.Range("C28").Value = 384
$ws.Range("D28").Value = 110324
$ws.Range("E28").Value = 3134

# Row 38: refreshed figures
$ws.Range("B38").Value = 81094
$ws.Range("C38").Value = 595
$ws.Range("D38").Value = 45666
$ws.Range("E38").Value = 34082
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 1346

# Row 47: refreshed figures
$ws.Range("D47").Value = 50128
$ws.Range("E47").Value = 5198

# Row 74: refreshed figures
$ws.Range("B74").Value = 18678
$ws.Range("C74").Value = 184
$ws.Range("D74").Value = 13214
$ws.Range("E74").Value = 5073
$ws.Range("G74").Value = 1
$ws.Range("H74").Value = 391

# Row 101: now "Grecia" (label + figures updated / row reordered)
$ws.Range("A101").Value = 'Grecia'
$ws.Range("B101").Value = 5942
$ws.Range("C101").Value = 193
$ws.Range("D101").Value = 3804
$ws.Range("E101").Value = 1924
$ws.Range("G101").Value = 1
$ws.Range("H101").Value = 214

# Row 102: now "Libia" (label + figures updated / row reordered)
$ws.Range("A102").Value = 'Libia'
$ws.Range("B102").Value = 5929
$ws.Range("D102").Value = 724
$ws.Range("E102").Value = 5080
$ws.Range("H102").Value = 125

# Row 122: refreshed figures
$ws.Range("B122").Value = 2880
$ws.Range("C122").Value = 9
$ws.Range("E122").Value = 247

# Row 144: refreshed figures
$ws.Range("B144").Value = 1283
$ws.Range("C144").Value = 15
$ws.Range("D144").Value = 1189
$ws.Range("E144").Value = 83

# Row 182: now "San Martin (Parte Holandesa)" (label + figures updated / row reordered)
$ws.Range("A182").Value = 'San Martin (Parte Holandesa)'
$ws.Range("B182").Value = 205
$ws.Range("C182").Value = 16
$ws.Range("D182").Value = 93
$ws.Range("E182").Value = 95
$ws.Range("H182").Value = 17

# Row 183: now "Islas Caimanes" (label + figures updated / row reordered)
$ws.Range("A183").Value = 'Islas Caimanes'
$ws.Range("B183").Value = 203
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 202
$ws.Range("E183").Value = 0
$ws.Range("H183").Value = 1

# Row 184: now "Gibraltar" (label + figures updated / row reordered)
$ws.Range("A184").Value = 'Gibraltar'
$ws.Range("B184").Value = 202
$ws.Range("C184").Value = 1
$ws.Range("D184").Value = 187
$ws.Range("E184").Value = 15
$ws.Range("H184").Value = 0

# Row 186: refreshed figures
$ws.Range("B186").Value = 159
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 6

# Row 187: refreshed figures
$ws.Range("B187").Value = 143
$ws.Range("C187").Value = 1
$ws.Range("E187").Value = 24

# Row 195: refreshed figures
$ws.Range("B195").Value = 84
$ws.Range("C195").Value = 6
$ws.Range("E195").Value = 37

# Row 213: now "Montserrat" (label + figures updated / row reordered)
$ws.Range("A213").Value = 'Montserrat'
$ws.Range("D213").Value = 12
$ws.Range("H213").Value = 1

# Row 214: now "Islas Malvinas" (label + figures updated / row reordered)
$ws.Range("A214").Value = 'Islas Malvinas'
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0
